$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "47.376.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.502.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.29%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.23%  "

$ws.Range("E7").Value = "  +1.24%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.23%  "

$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.894.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.512.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.861"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.319.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.77%  "

$ws.Range("E20").Value = "  +4.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0942"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "248.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.26%  "

$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("E28").Value = "  +0.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.139"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.37%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0796"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.89%  "

$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.07%  "

$ws.Range("E39").Value = "  +1.88%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.112"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.57%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "121.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.85%  "

$ws.Range("E42").Value = "  -1.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.22%  "

$ws.Range("E44").Value = "  +2.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.999.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "

$ws.Range("E48").Value = "  -4.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.33%  "

$ws.Range("E50").Value = "  +5.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.08%  "
